# "Added Alembic and modified ppt"
# Slide 3 ("Technology Stack:") lists the backend tech stack as one
# paragraph per line inside the "Content Placeholder 2" shape:
#   FastAPI / PostgreSQL / SQLAlchemy / Alembic / Pydantic / JWT / SMTP
# The edit removes the standalone "Alembic" line (its info is folded
# back in / no longer called out separately), leaving:
#   FastAPI / PostgreSQL / SQLAlchemy / Pydantic / JWT / SMTP

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Locate the "Alembic" paragraph robustly (rather than hard-coding an
# index) and remove it entirely, which also merges it away cleanly.
for ($i = $tr.Paragraphs().Count; $i -ge 1; $i--) {
    $para = $tr.Paragraphs($i)
    if ($para.Text.Trim() -eq "Alembic") {
        $para.Delete()
    }
}
